$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New diary entry for "6 marras" (row 25) ------------------------------
# Write the date first so its shared string lands before the corrected
# "Mathematics of Rotations" note below.
$ws.Range("A25").Value = "6 marras"

# Row 24 ("5 marras"): the reading note is corrected - the chapter was only
# read up to page 161 that day (the rest moved into the new 6 marras entry).
$ws.Range("C24").Value = '"The Mathematics of Rotations", oppikirjasta 145-161'

# Continuing content for the new row: more matrix operations implemented -
# inverse, transpose, determinant for 3x4.
$ws.Range("C25").Value = "oppikirjasta 162- , Implementing mathematics of rotations, "
$ws.Range("B25").Value = "9.15-11.15, 12:15"
$ws.Range("G25").Value = 3.5

# Match the formatting used by the other rows in the diary table: time
# format + wrap-text for "Kello" (col B), wrap-text for the content columns.
$ws.Range("B25").NumberFormat = $ws.Range("B20").NumberFormat
$ws.Range("B25").WrapText = $true
$ws.Range("C25").WrapText = $true

$ws.Rows("25").RowHeight = 29

# Move the selection/view the way Excel would after typing the new row.
[void]$ws.Range("D25").Select()

$wb.Save()
